$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 3
    4  = 0
    5  = 0
    6  = 0
    7  = 2
    8  = 3
    9  = 0
    10 = 3
    11 = 2
    12 = 0
    13 = 0
    14 = 2
    15 = 3
    16 = 1
    17 = 4
    18 = 3
    19 = 2
    20 = 1
    21 = 1
    22 = 4
    23 = 0
    24 = 0
    25 = 1
    26 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
